# Rename the sheet from "Sheet1" to "Estimates" and fix the label in A65
# from "Min (P=95%)" (duplicate) to "Max (P=95%)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Estimates"
$ws.Range("A65").Value = "Max (P=95%)"
